# Update countries & provincias Spain
# Applies the COVID data refresh: updated case numbers for several
# countries, which causes Italia to overtake Irak (rows 22/23) and
# Uganda to overtake Republica de Africa Central & Cabo Verde
# (rows 119/120/121); also bumps the "Datos actualizados" footer time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 6683812
$ws.Range("C4").Value = 7211
$ws.Range("D4").Value = 3950642
$ws.Range("E4").Value = 2534990
$ws.Range("G4").Value = 52
$ws.Range("H4").Value = 198180

# --- Row 5: India ---
$ws.Range("B5").Value = 4788593
$ws.Range("C5").Value = 36805
$ws.Range("D5").Value = 3730949
$ws.Range("E5").Value = 978713
$ws.Range("G5").Value = 317
$ws.Range("H5").Value = 78931

# --- Rows 22/23: Italia overtakes Irak ---
$ws.Range("A22").Value = "Italia"
$ws.Range("B22").Value = 287753
$ws.Range("C22").Value = 1456
$ws.Range("D22").Value = 213634
$ws.Range("E22").Value = 38509
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 35610

$ws.Range("A23").Value = "Irak"
$ws.Range("B23").Value = 286778
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 221283
$ws.Range("E23").Value = 57554
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 7941

# --- Row 25 ---
$ws.Range("B25").Value = 260787
$ws.Range("C25").Value = 241
$ws.Range("E25").Value = 16060

# --- Row 29 ---
$ws.Range("B29").Value = 136345
$ws.Range("C29").Value = 204
$ws.Range("D29").Value = 120199
$ws.Range("E29").Value = 6975
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 9171

# --- Row 55 ---
$ws.Range("D55").Value = 56764
$ws.Range("E55").Value = 615

# --- Row 65 ---
$ws.Range("B65").Value = 42978
$ws.Range("C65").Value = 264
$ws.Range("E65").Value = 11418
$ws.Range("G65").Value = 6
$ws.Range("H65").Value = 1123

# --- Row 67 ---
$ws.Range("B67").Value = 38327
$ws.Range("C67").Value = 155
$ws.Range("D67").Value = 35756
$ws.Range("E67").Value = 2009
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 562

# --- Row 93 ---
$ws.Range("B93").Value = 12123
$ws.Range("C93").Value = 44
$ws.Range("E93").Value = 1487

# --- Row 101 ---
$ws.Range("B101").Value = 9049
$ws.Range("C101").Value = 35
$ws.Range("D101").Value = 7816
$ws.Range("E101").Value = 1161

# --- Row 112 ---
$ws.Range("B112").Value = 5395
$ws.Range("C112").Value = 1
$ws.Range("D112").Value = 5330
$ws.Range("E112").Value = 4

# --- Rows 119/120/121: Uganda overtakes Republica de Africa Central & Cabo Verde ---
$ws.Range("A119").Value = "Uganda"
$ws.Range("B119").Value = 4799
$ws.Range("C119").Value = 96
$ws.Range("D119").Value = 2256
$ws.Range("E119").Value = 2488
$ws.Range("G119").Value = 3
$ws.Range("H119").Value = 55

$ws.Range("A120").Value = "Republica de Africa Central"
$ws.Range("B120").Value = 4749
$ws.Range("D120").Value = 1825
$ws.Range("E120").Value = 2862
$ws.Range("H120").Value = 62

$ws.Range("A121").Value = "Cabo Verde"
$ws.Range("B121").Value = 4711
$ws.Range("D121").Value = 4104
$ws.Range("E121").Value = 563
$ws.Range("H121").Value = 44

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 17:17"
